$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived; insert it at row 10 (pushing the
# existing data rows 10-72 down to 11-73, same as the prior week's rows
# that already occupy rows 10 onward).
$ws.Rows.Item(10).Insert()

# Fill in the new row 10 with the new data (a new weekly price entry)
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(10, 3).Value = "Bíobío"
$ws.Cells.Item(10, 4).Value = 45069
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 100112043
$ws.Cells.Item(10, 7).Value = "Pepino dulce"
$ws.Cells.Item(10, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 150
$ws.Cells.Item(10, 11).Value = 16000
$ws.Cells.Item(10, 12).Value = 17000
$ws.Cells.Item(10, 13).Value = 16333
$ws.Cells.Item(10, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 907
$ws.Cells.Item(10, 17).Value = 18
$ws.Cells.Item(10, 18).Value = "Hortaliza"
